$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1382.6666
$ws.Range("J58").Value = 1450
$ws.Range("L58").Value = 4350
$ws.Range("N58").Value = -4650
$ws.Range("H70").Value = 6527.091
$ws.Range("J70").Value = 6849.75
$ws.Range("L70").Value = 20549.25
$ws.Range("N70").Value = -21089.25
$ws.Range("H73").Value = 6527.091
$ws.Range("J73").Value = 6849.75
$ws.Range("L73").Value = 20549.25
$ws.Range("N73").Value = -22421.25
$ws.Range("H86").Value = 4250
$ws.Range("I86").Value = 4333.3335
$ws.Range("K86").Value = 4333.3335
$ws.Range("M86").Value = -3210.3335
$ws.Range("H89").Value = 4250
$ws.Range("I89").Value = 4333.3335
$ws.Range("K89").Value = 21666.6675
$ws.Range("M89").Value = -16050.6675
$ws.Range("H116").Value = 4233.3335
$ws.Range("I116").Value = 4225
$ws.Range("K116").Value = 4225
$ws.Range("M116").Value = -783
$ws.Range("H132").Value = 1195.2858
$ws.Range("I132").Value = 1184
$ws.Range("K132").Value = 3552
$ws.Range("M132").Value = -1022
$ws.Range("H133").Value = 119709
$ws.Range("I133").Value = 119709
$ws.Range("K133").Value = 119709
$ws.Range("M133").Value = -114649
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H138").Value = 2615.6978
$ws.Range("I138").Value = 1793.3572
$ws.Range("J138").Value = 3012.6897
$ws.Range("K138").Value = 5380.071599999999
$ws.Range("L138").Value = 9038.069100000001
$ws.Range("M138").Value = -240.0715999999993
$ws.Range("N138").Value = -19318.0691
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 395.7143
$ws.Range("I5").Value = 239.33333
$ws.Range("J5").Value = 677.2
$ws.Range("K5").Value = 239.33333
$ws.Range("L5").Value = 677.2
$ws.Range("M5").Value = -127.33333
$ws.Range("N5").Value = -901.2
$ws.Range("H63").Value = 7922.3213
$ws.Range("I63").Value = 2574.2856
$ws.Range("K63").Value = 2574.2856
$ws.Range("M63").Value = -1888.2856
$ws.Range("H66").Value = 7922.3213
$ws.Range("I66").Value = 2574.2856
$ws.Range("K66").Value = 12871.428
$ws.Range("M66").Value = -9439.428
$ws.Range("H75").Value = 55277
$ws.Range("J75").Value = 55277
$ws.Range("L75").Value = 55277
$ws.Range("N75").Value = -57025
$ws.Range("H78").Value = 55277
$ws.Range("J78").Value = 55277
$ws.Range("L78").Value = 165831
$ws.Range("N78").Value = -174567
$ws.Range("H122").Value = 5882
$ws.Range("I122").Value = 4500
$ws.Range("K122").Value = 13500
$ws.Range("M122").Value = -11050
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 395.7143
$ws.Range("I4").Value = 239.33333
$ws.Range("J4").Value = 677.2
$ws.Range("K4").Value = 239.33333
$ws.Range("L4").Value = 677.2
$ws.Range("M4").Value = -124.33333
$ws.Range("N4").Value = -907.2
$ws.Range("H20").Value = 2475.9644
$ws.Range("I20").Value = 2319.3333
$ws.Range("J20").Value = 2656.6924
$ws.Range("K20").Value = 2319.3333
$ws.Range("L20").Value = 2656.6924
$ws.Range("M20").Value = -2072.3333
$ws.Range("N20").Value = -3150.6924
$ws.Range("H21").Value = 45135.5
$ws.Range("J21").Value = 45135.5
$ws.Range("L21").Value = 45135.5
$ws.Range("N21").Value = -45607.5
$ws.Range("H132").Value = 119995
$ws.Range("J132").Value = 119995
$ws.Range("L132").Value = 119995
$ws.Range("N132").Value = -130115
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 33905
$ws.Range("J28").Value = 33905
$ws.Range("L28").Value = 33905
$ws.Range("N28").Value = -34395
$ws.Range("H62").Value = 9049.5
$ws.Range("I62").Value = 4333
$ws.Range("K62").Value = 4333
$ws.Range("M62").Value = -3709
$ws.Range("H65").Value = 9049.5
$ws.Range("I65").Value = 4333
$ws.Range("K65").Value = 21665
$ws.Range("M65").Value = -18545
$ws.Range("H123").Value = 97000
$ws.Range("J123").Value = 97000
$ws.Range("L123").Value = 97000
$ws.Range("N123").Value = -106800
$ws.Range("H134").Value = 8870.412
$ws.Range("I134").Value = 9095.643
$ws.Range("J134").Value = 7819.3335
$ws.Range("K134").Value = 27286.929
$ws.Range("L134").Value = 23458.0005
$ws.Range("M134").Value = -24751.929
$ws.Range("N134").Value = -28528.0005
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2914.8
$ws.Range("J12").Value = 2683.3333
$ws.Range("L12").Value = 8049.999899999999
$ws.Range("N12").Value = -8395.999899999999
$ws.Range("H15").Value = 3697.516
$ws.Range("I15").Value = 492.7
$ws.Range("J15").Value = 5223.619
$ws.Range("K15").Value = 1478.1
$ws.Range("L15").Value = 15670.857
$ws.Range("M15").Value = -1338.1
$ws.Range("N15").Value = -15950.857
$ws.Range("H64").Value = 1000
$ws.Range("I64").Value = 1000
$ws.Range("K64").Value = 3000
$ws.Range("M64").Value = -2730
$ws.Range("H67").Value = 1000
$ws.Range("I67").Value = 1000
$ws.Range("K67").Value = 3000
$ws.Range("M67").Value = -2064
$ws.Range("H74").Value = 12012.6
$ws.Range("J74").Value = 12262.5
$ws.Range("L74").Value = 36787.5
$ws.Range("N74").Value = -38909.5
$ws.Range("H77").Value = 12012.6
$ws.Range("J77").Value = 12262.5
$ws.Range("L77").Value = 110362.5
$ws.Range("N77").Value = -120970.5
$ws.Range("H81").Value = 5933.3335
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 7900
$ws.Range("K81").Value = 6000
$ws.Range("L81").Value = 23700
$ws.Range("M81").Value = -4877
$ws.Range("N81").Value = -25946
$ws.Range("H84").Value = 5933.3335
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 7900
$ws.Range("K84").Value = 18000
$ws.Range("L84").Value = 71100
$ws.Range("M84").Value = -12384
$ws.Range("N84").Value = -82332
$ws.Range("H87").Value = 500
$ws.Range("I87").Value = 500
$ws.Range("K87").Value = 1500
$ws.Range("M87").Value = -252
$ws.Range("H90").Value = 500
$ws.Range("I90").Value = 500
$ws.Range("K90").Value = 4500
$ws.Range("M90").Value = 1740
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 38941.54
$ws.Range("J2").Value = 111239.78
$ws.Range("L2").Value = 111239.78
$ws.Range("N2").Value = -111465.78
$ws.Range("H48").Value = 14996.667
$ws.Range("J48").Value = 14996.667
$ws.Range("L48").Value = 14996.667
$ws.Range("N48").Value = -15966.667
$ws.Range("H49").Value = 59995
$ws.Range("J49").Value = 59995
$ws.Range("L49").Value = 59995
$ws.Range("N49").Value = -60363
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H102").Value = 6170.6665
$ws.Range("I102").Value = 6170.6665
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 6170.6665
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -4548.6665
$ws.Range("N102").ClearContents()
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H122").Value = 11210
$ws.Range("I122").Value = 11859.111
$ws.Range("K122").Value = 35577.333
$ws.Range("M122").Value = -33127.333
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 16682333
$ws.Range("I42").Value = 50000000
$ws.Range("J42").Value = 23500
$ws.Range("K42").Value = 50000000
$ws.Range("L42").Value = 23500
$ws.Range("M42").Value = -49999437
$ws.Range("N42").Value = -24626
$ws.Range("H49").Value = 16682333
$ws.Range("I49").Value = 50000000
$ws.Range("J49").Value = 23500
$ws.Range("K49").Value = 50000000
$ws.Range("L49").Value = 23500
$ws.Range("M49").Value = -49999853
$ws.Range("N49").Value = -23794
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H112").Value = 84187
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 84187
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 84187
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -87141
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 45000
$ws.Range("J56").Value = 45000
$ws.Range("L56").Value = 45000
$ws.Range("N56").Value = -46428
$ws.Range("H122").Value = 7083.6665
$ws.Range("I122").Value = 7083.6665
$ws.Range("K122").Value = 21250.9995
$ws.Range("M122").Value = -18800.9995
$ws.Range("H126").Value = 106104.3
$ws.Range("J126").Value = 7473.5
$ws.Range("L126").Value = 22420.5
$ws.Range("N126").Value = -27360.5
$ws.Range("H132").Value = 6266.357
$ws.Range("I132").Value = 2505.4614
$ws.Range("K132").Value = 7516.3842
$ws.Range("M132").Value = -4986.3842

Write-Output "Applied 236 cell updates across 8 sheets"